# Update District column (G) values to the cleaned "Ballari (Bellary)" label
# for the rows whose value was still "Ballari", "Ballary", or an incorrect
# value (school name instead of district) left over from prior cleanup.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(3,9,11,13,14,16,17,18,20,21,26,28,29,33,34,35,36,41,42,45,46,48,49,50,51,52,57,58,59)

foreach ($r in $rows) {
    $ws.Range("G$r").Value = "Ballari (Bellary)"
}
